$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newText = "Test - Test - 10/30/2020"

foreach ($addr in @("B2", "B4", "B8", "B12", "B16")) {
    $ws.Range($addr).Value = $newText
}

$ws.Range("B16").Select()
